# Update the cryptos list with the latest scraped prices / volumes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text (e.g. "313.66", "1.002").
# Force the whole column to text format BEFORE assigning new values so Excel
# doesn't silently convert these strings into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "28.135.32"
$ws.Cells.Item(2, 5).Value = "  -1.43%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.796.44"
$ws.Cells.Item(3, 5).Value = "  -1.56%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.12%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "313.66"
$ws.Cells.Item(5, 5).Value = "  -0.59%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.12%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5216"
$ws.Cells.Item(7, 5).Value = "  +2.00%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.3835"
$ws.Cells.Item(8, 5).Value = "  -3.01%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.07896"
$ws.Cells.Item(9, 5).Value = "  -3.96%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "41.33"

# Row 11
$ws.Cells.Item(11, 4).Value = "1.101"
$ws.Cells.Item(11, 5).Value = "  -1.37%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "6.288"
$ws.Cells.Item(12, 5).Value = "  -0.88%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.002"
$ws.Cells.Item(13, 5).Value = "  +0.12%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "20.57"
$ws.Cells.Item(14, 5).Value = "  -2.91%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "1.794.07"
$ws.Cells.Item(15, 5).Value = "  -1.69%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "7.286"
$ws.Cells.Item(16, 5).Value = "  -3.57%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "92.77"
$ws.Cells.Item(17, 5).Value = "  -0.19%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "0.00001086"
$ws.Cells.Item(18, 5).Value = "  -3.51%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "0.06550"
$ws.Cells.Item(19, 5).Value = "  -1.48%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "1.002"
$ws.Cells.Item(20, 5).Value = "  +0.16%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "17.32"
$ws.Cells.Item(21, 5).Value = "  -2.97%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "5.959"
$ws.Cells.Item(22, 5).Value = "  -2.25%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "28.159.21"
$ws.Cells.Item(23, 5).Value = "  -1.48%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "11.14"
$ws.Cells.Item(24, 5).Value = "  -2.65%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.26%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "161.12"
$ws.Cells.Item(26, 5).Value = "  +2.95%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "20.51"
$ws.Cells.Item(27, 5).Value = "  -4.07%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "1.998.28"
$ws.Cells.Item(28, 5).Value = "  -1.82%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "2.334"
$ws.Cells.Item(29, 5).Value = "  -3.49%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "122.96"
$ws.Cells.Item(30, 5).Value = "  -3.14%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.1065"
$ws.Cells.Item(31, 5).Value = "  -2.12%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "1.054"
$ws.Cells.Item(32, 5).Value = "  -5.43%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "3.674"
$ws.Cells.Item(33, 5).Value = "  +0.45%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "5.550"
$ws.Cells.Item(34, 5).Value = "  -3.79%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.07280"
$ws.Cells.Item(35, 5).Value = "  +3.28%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "12.26"
$ws.Cells.Item(36, 5).Value = "  +8.59%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "0.02320"
$ws.Cells.Item(37, 5).Value = "  -1.54%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "FraxShare"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(38, 4).Value = "8.766"
$ws.Cells.Item(38, 5).Value = "  -0.49%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Algorand"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(39, 4).Value = "0.2140"
$ws.Cells.Item(39, 5).Value = "  -4.09%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "5.078"
$ws.Cells.Item(40, 5).Value = "  -3.91%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "0.6158"
$ws.Cells.Item(41, 5).Value = "  -2.79%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "1.161"
$ws.Cells.Item(42, 5).Value = "  -1.89%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.378"
$ws.Cells.Item(43, 5).Value = "  -1.54%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "13.24"
$ws.Cells.Item(44, 5).Value = "  -2.30%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "3.773"
$ws.Cells.Item(45, 5).Value = "  +1.04%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "0.5952"
$ws.Cells.Item(46, 5).Value = "  +0.04%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "127.64"
$ws.Cells.Item(47, 5).Value = "  +1.99%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "1.231"
$ws.Cells.Item(48, 5).Value = "  +3.26%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "1.919"
$ws.Cells.Item(49, 5).Value = "  -3.93%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "0.06753"
$ws.Cells.Item(50, 5).Value = "  -2.81%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "73.13"
$ws.Cells.Item(51, 5).Value = "  -1.59%  "
